$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "54.776.93"
$ws.Range("E2").Value = "  +1.18%  "

$ws.Range("D3").Value = "2.272.43"
$ws.Range("E3").Value = "  +0.12%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "503.01"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.06%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "127.78"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.31%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.20%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.528"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.42%  "

$ws.Range("D9").Value = "2.289.46"
$ws.Range("E9").Value = "  +0.77%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0973"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +2.38%  "

$ws.Range("E11").Value = "  +1.17%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.09"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +8.09%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.338"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.24%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.37"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +3.18%  "

$ws.Range("D15").Value = "2.678.35"
$ws.Range("E15").Value = "  +0.28%  "

$ws.Range("D16").Value = "54.890.98"
$ws.Range("E16").Value = "  +1.43%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000130"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.58%  "

$ws.Range("D18").Value = "2.277.37"
$ws.Range("E18").Value = "  +0.49%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.35"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.42%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.18"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.08%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "309.44"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +2.08%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.54"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +3.61%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.998"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.28%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "59.89"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.92%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.995"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.44%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.154"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +2.93%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.45"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +2.28%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "172.27"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.92%  "

$ws.Range("B29").Value = "Aptos"
$ws.Range("C29").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.09"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +2.91%  "

$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.63"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.34%  "

$ws.Range("B31").Value = "PEPE"
$ws.Range("C31").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D31").Value = "0.0₃0699"
$ws.Range("E31").Value = "  +1.43%  "

$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.14"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +5.60%  "

$ws.Range("E33").Value = "  +0.00%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "17.85"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.72%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.997"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.16%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.22"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +2.48%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.899"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -4.91%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.85"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +4.38%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.80"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +2.21%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.44"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +3.10%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.375"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.16%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "135.36"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +8.48%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.44"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +2.72%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.84"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.63%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "256.92"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +7.63%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0503"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +2.59%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0911"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +2.38%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.546"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.28%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.374"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +1.04%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0211"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +2.90%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "16.39"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.55%  "
